$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix K268: -1 -> 0
$ws.Range("K268").Value = 0

# Duplicate row 268 (template row, now with corrected K=0 and full formatting/style)
# into new rows 269-289, one insert at a time so each new row keeps the same styling
# (s="6" on column E, etc.) as the source row.
for ($i = 0; $i -lt 21; $i++) {
    $ws.Rows("268:268").Copy()
    $ws.Rows("269:269").Insert()
}

# Now overwrite only the cells that differ from the template row 268 on each new row
# Row 276
$ws.Range("D276").Value = "CE"

# Row 277
$ws.Range("A277").Value = 39400
$ws.Range("C277").Value = 39470
$ws.Range("F277").Value = 39530
$ws.Range("G277").Value = 39420

# Row 278
$ws.Range("A278").Value = 39400
$ws.Range("C278").Value = 39470
$ws.Range("F278").Value = 39530
$ws.Range("G278").Value = 39420

# Row 279
$ws.Range("A279").Value = 39400
$ws.Range("C279").Value = 39435
$ws.Range("G279").Value = 39380

# Row 280
$ws.Range("A280").Value = 39400
$ws.Range("C280").Value = 39425
$ws.Range("G280").Value = 39380

# Row 281
$ws.Range("A281").Value = 39400
$ws.Range("C281").Value = 39435
$ws.Range("G281").Value = 39380

# Row 282
$ws.Range("A282").Value = 39400
$ws.Range("C282").Value = 39435
$ws.Range("G282").Value = 39380

# Row 283
$ws.Range("A283").Value = 38400
$ws.Range("C283").Value = 38455
$ws.Range("G283").Value = 39380

# Row 284
$ws.Range("A284").Value = 38400
$ws.Range("C284").Value = 38470
$ws.Range("F284").Value = 38500
$ws.Range("G284").Value = 38430

# Row 285
$ws.Range("A285").Value = 38500
$ws.Range("C285").Value = 38470
$ws.Range("D285").Value = "CE"
$ws.Range("F285").Value = 38430
$ws.Range("G285").Value = 38500

# Row 286
$ws.Range("A286").Value = 38500
$ws.Range("C286").Value = 38470
$ws.Range("D286").Value = "CE"
$ws.Range("F286").Value = 38430
$ws.Range("G286").Value = 38500

# Row 287
$ws.Range("A287").Value = 38450
$ws.Range("C287").Value = 38470
$ws.Range("F287").Value = 38600
$ws.Range("G287").Value = 38350

# Row 288
$ws.Range("A288").Value = 38450
$ws.Range("C288").Value = 38470
$ws.Range("F288").Value = 38600
$ws.Range("G288").Value = 38350

# Row 289
$ws.Range("A289").Value = 38450
$ws.Range("C289").Value = 38470
$ws.Range("F289").Value = 38600
$ws.Range("G289").Value = 38350

